$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 31 cells (added "2c" -> 2*F24 line for two-chord adimensionalization) ---
$ws.Range("E31").Value2 = "2c"
$ws.Range("F31").Formula = "=2*F24"
$ws.Range("G31").Value2 = "m"

# --- Unit labels "m" added next to F24/F25/F28 ---
$ws.Range("G24").Value2 = "m"
$ws.Range("G25").Value2 = "m"
$ws.Range("G28").Value2 = "m"

# --- Update formulas to reference F31 (2*F24) instead of F24 ---
$ws.Range("J24").Formula = "=(F27 + (F28/F31)*(1-F29))/(1+(1-F29))"
$ws.Range("J25").Formula = "=F25+J24*F31"

# --- Sheet view change: scroll/selection moved ---
$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 7
$activeWindow.ScrollColumn = 1
$ws.Range("J27").Select()
